# The "congenital" tag used in the per-variable category lists (column A)
# has been renamed to "misc_long_term" across the whole workbook.
$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rows = $used.Rows.Count
    for ($r = 1; $r -le $rows; $r++) {
        $cell = $ws.Cells.Item($r, 1)
        if ($cell.Value() -eq "congenital") {
            $cell.Value = "misc_long_term"
        }
    }
}
